$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dallin")

# --- Header row additions (F1, G1) ---
$ws.Range("F1").Value = "Start `n(Center Column)"
$ws.Range("G1").Value = "End `n(Center Column)"
$ws.Range("F1:G1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30.5
$ws.Columns.Item(6).ColumnWidth = 14.3
$ws.Columns.Item(7).ColumnWidth = 14.6

# --- Insert two new blank rows: before old row 11, and before old row 12 (now 14) ---
# After first insert, old row 11 (D11/E11=1301/1372) becomes row 12.
$ws.Rows.Item(11).Insert()
# After the first insert, old row 12 (D12/E12=1379/1410) is now row 13; insert before it too.
$ws.Rows.Item(13).Insert()

# --- Orbit 2 (row 2-3) unchanged data, just B2 text stays the same string ---
$ws.Range("B2").Value = "awe_l1r_q20_2023326T0108_00002_v01.nc"

# --- Orbit 4 (rows 4-7): add centerline columns F/G on rows 4,5,6 ---
$ws.Range("B4").Value = "awe_l1r_q20_2023326T0108_00004_v01.nc"
$ws.Range("F4").Value = 912
$ws.Range("G4").Value = 922
$ws.Range("F5").Value = 941
$ws.Range("G5").Value = 1057
$ws.Range("F6").Value = 1194
$ws.Range("G6").Value = 1213

# --- Orbit 6 (row 8): add F8/G8 ---
$ws.Range("B8").Value = "awe_l1r_q20_2023326T0108_00006_v01.nc"
$ws.Range("F8").Value = 987
$ws.Range("G8").Value = 1078

# --- Orbit 8 (rows 9-14): add F/G data on rows 10-14 (new rows 11 & 13 included) ---
$ws.Range("B9").Value = "awe_l1r_q20_2023326T0108_00008_v01.nc"
$ws.Range("F10").Value = 1041
$ws.Range("G10").Value = 1049
$ws.Range("F11").Value = 1060
$ws.Range("G11").Value = 1150
$ws.Range("F12").Value = 1313
$ws.Range("G12").Value = 1315
$ws.Range("F13").Value = 1325
$ws.Range("G13").Value = 1330
$ws.Range("F14").Value = 1410
$ws.Range("G14").Value = 1410

# --- Orbit 10 (row 15-16): B text only, no centerline data ---
$ws.Range("B15").Value = "awe_l1r_q20_2023326T0108_00010_v01.nc"

# --- Orbit 12 (row 17) ---
$ws.Range("B17").Value = "awe_l1r_q20_2023326T0108_00012_v01.nc"

# --- Orbit 14 (rows 18-19) ---
$ws.Range("B18").Value = "awe_l1r_q20_2023326T0108_00014_v01.nc"

# --- Orbit 15 (rows 20-21) ---
$ws.Range("B20").Value = "awe_l1r_q20_2023326T0108_00015_v01.nc"

# --- New: Orbit 16 (rows 22-26) ---
$ws.Range("A22").Value = 16
$ws.Range("B22").Value = "awe_l1r_q20_2023326T0108_00016_v01.nc"
$ws.Range("C22").Value = 1406
$ws.Range("D22").Value = 982
$ws.Range("E22").Value = 1406
$ws.Range("F22").Value = 1015
$ws.Range("G22").Value = 1043

$ws.Range("F23").Value = 1092
$ws.Range("G23").Value = 1111

$ws.Range("F24").Value = 1160
$ws.Range("G24").Value = 1233

$ws.Range("F25").Value = 1270
$ws.Range("G25").Value = 1321

$ws.Range("F26").Value = 1331
$ws.Range("G26").Value = 1390

# --- New: Orbit 18 (row 27) ---
$ws.Range("A27").Value = 18
$ws.Range("B27").Value = "awe_l1r_q20_2023326T0108_00018_v01.nc"
$ws.Range("C27").Value = 1414
$ws.Range("D27").Value = 1414
$ws.Range("E27").Value = 1122
$ws.Range("F27").Value = 1122
$ws.Range("G27").Value = 1414

# --- New: Orbit 20 (rows 28-35) ---
$ws.Range("A28").Value = 20
$ws.Range("B28").Value = "awe_l1r_q20_2023326T0108_00020_v01.nc"
$ws.Range("C28").Value = 1419
$ws.Range("D28").Value = 782
$ws.Range("E28").Value = 900
$ws.Range("F28").Value = 819
$ws.Range("G28").Value = 859

$ws.Range("D29").Value = 979
$ws.Range("E29").Value = 1192
$ws.Range("F29").Value = 1018
$ws.Range("G29").Value = 1020

$ws.Range("F30").Value = 1038
$ws.Range("G30").Value = 1047

$ws.Range("F31").Value = 1053
$ws.Range("G31").Value = 1054

$ws.Range("F32").Value = 1074
$ws.Range("G32").Value = 1128

$ws.Range("F33").Value = 1145
$ws.Range("G33").Value = 1149

$ws.Range("D34").Value = 1288
$ws.Range("E34").Value = 1419
$ws.Range("F34").Value = 1320
$ws.Range("G34").Value = 1323

$ws.Range("F35").Value = 1332
$ws.Range("G35").Value = 1419

# --- Final selection state ---
$ws.Range("C39").Select()
